$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The BF column holds the game date as a literal text string (e.g. "2013-05-23").
# Writing that string straight into .Value lets Excel auto-detect it as a date
# and convert it to a date serial, which is not what we want here - the source
# data is a plain "YYYY-MM-DD" label, not a real date value. Force the cell to
# text first, write the corrected label, then restore the cell's style so we
# don't leave a stray number-format behind on the cell.
for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    $cell.NumberFormat = "@"
    $cell.Value = "2013-05-23"
    $cell.Style = "Normal"
}
